$d = $word.ActiveDocument

# The second paragraph of the document currently reads:
#   "N+1 problem is a performance issue in Object Relational Mapping that
#    fires multiple select queries ..."
# It must be replaced with new wording, rendered in a different font /
# size / color / shading, and followed by a trailing run (a single space)
# that keeps the paragraph's original bold+underlined sz32 look (which is
# also what the paragraph mark's own formatting should end up as).
#
# We rebuild the paragraph's XML precisely (preserving its original
# paraId/rsid attributes) and hand it to Range.InsertXML, which replaces
# the contents of the supplied Range with the given WordprocessingML.

$targetParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4F3586C0" w14:textId="0477A079" w:rsidR="00D34067" w:rsidRPr="007F31F0" w:rsidRDefault="00D34067" w:rsidP="007F31F0"><w:pPr><w:pStyle w:val="NoSpacing"/><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/><w:color w:val="242424"/><w:spacing w:val="-1"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>The N+1 loading problem occurs when an application makes one query to retrieve the initial data (e.g., a list of users) and then an additional query for each row of data to retrieve related data (e.g., the user’s articles) where those related data were not part of the initial data. Hence, if there are N users, the application will make 1 query to get all users and then N additional queries to get the articles for each user, totaling N+1 queries.</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@

$targetParagraph = $d.Paragraphs.Item(2)
$targetParagraph.Range.InsertXML($targetParaXml)

Write-Output $d.Paragraphs.Item(2).Range.Text
